# Update the coupon code for MWL (row 2, column A / "id") on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = "CA-XGYTNHX3"
